# Updates cryptos list price (D) / 1h volume-change (E) columns, and three
# coin rows that moved position (B/C/D/E), per the Feb 16 2024 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (matches the source file, where every cell
# in these columns is stored as a string, not a number). A bare numeric-
# looking string (e.g. "40.27") would otherwise be auto-converted to a
# number by Excel, so we use the classic leading quote-prefix, then reset
# the style back to Normal so no stray number-format/quote-prefix survives.
function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '52.134.06'
$ws.Range('E2').Value = '  +0.46%  '
$ws.Range('D3').Value = '2.800.37'
$ws.Range('E3').Value = '  -0.84%  '
$ws.Range('E4').Value = '  -0.03%  '
Set-TextValue $ws.Range('D5') '360.39'
Set-TextValue $ws.Range('D6') '110.77'
$ws.Range('E6').Value = '  -2.34%  '
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('E8').Value = '  +0.00%  '
Set-TextValue $ws.Range('D9') '0.597'
$ws.Range('E9').Value = '  -1.23%  '
Set-TextValue $ws.Range('D10') '40.27'
$ws.Range('E10').Value = '  -3.21%  '
$ws.Range('E12').Value = '  +1.34%  '
$ws.Range('E13').Value = '  -1.97%  '
Set-TextValue $ws.Range('D14') '7.64'
$ws.Range('E14').Value = '  -1.65%  '
$ws.Range('D15').Value = '3.234.19'
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range('D16') '0.953'
$ws.Range('E16').Value = '  +7.22%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '2.765.05'
$ws.Range('E17').Value = '  -2.65%  '
$ws.Range('D18').Value = '52.006.70'
$ws.Range('E18').Value = '  +0.37%  '
Set-TextValue $ws.Range('D19') '7.44'
$ws.Range('E19').Value = '  -0.77%  '
$ws.Range('E20').Value = '  -1.38%  '
Set-TextValue $ws.Range('D21') '13.10'
$ws.Range('E21').Value = '  -2.50%  '
$ws.Range('D22').Value = '0.0₃0989'
$ws.Range('E22').Value = '  -0.02%  '
Set-TextValue $ws.Range('D23') '274.59'
Set-TextValue $ws.Range('D24') '70.40'
$ws.Range('E24').Value = '  +0.89%  '
$ws.Range('E25').Value = '  -1.11%  '
Set-TextValue $ws.Range('D26') '26.78'
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range('D28') '10.23'
$ws.Range('E28').Value = '  -0.97%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D29') '0.147'
$ws.Range('E29').Value = '  +5.38%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range('D30') '2.22'
$ws.Range('E30').Value = '  -1.38%  '
Set-TextValue $ws.Range('D31') '51.69'
$ws.Range('E31').Value = '  +1.84%  '
$ws.Range('E32').Value = '  +1.53%  '
Set-TextValue $ws.Range('D33') '34.53'
$ws.Range('E33').Value = '  +1.46%  '
Set-TextValue $ws.Range('D34') '5.79'
$ws.Range('E34').Value = '  -1.41%  '
Set-TextValue $ws.Range('D35') '0.0853'
$ws.Range('E35').Value = '  +3.26%  '
Set-TextValue $ws.Range('D36') '5.29'
$ws.Range('E36').Value = '  -0.49%  '
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('E38').Value = '  +0.35%  '
Set-TextValue $ws.Range('D39') '18.45'
$ws.Range('E39').Value = '  +1.06%  '
$ws.Range('E40').Value = '  -2.88%  '
Set-TextValue $ws.Range('D41') '2.60'
$ws.Range('E41').Value = '  +3.00%  '
$ws.Range('E42').Value = '  -1.15%  '
Set-TextValue $ws.Range('D43') '123.09'
$ws.Range('E43').Value = '  -2.66%  '
$ws.Range('E44').Value = '  -2.33%  '
Set-TextValue $ws.Range('D45') '22.30'
$ws.Range('E45').Value = '  -6.87%  '
$ws.Range('D46').Value = '2.082.44'
$ws.Range('E46').Value = '  -0.58%  '
Set-TextValue $ws.Range('D47') '3.28'
$ws.Range('E47').Value = '  -2.61%  '
$ws.Range('E48').Value = '  -2.38%  '
$ws.Range('E49').Value = '  +0.78%  '
Set-TextValue $ws.Range('D50') '0.935'
$ws.Range('E50').Value = '  -0.16%  '
$ws.Range('E51').Value = '  +0.37%  '
